$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Slit2"
$ws.Range("C2").Value = "Robo4"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.1645843333333333
$ws.Range("H2").Value = 0.493753
$ws.Range("I2").Value = 0.03485847193389392
$ws.Range("J2").Value = 0.03485847193389392
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 43.55927533333334
$ws.Range("N2").Value = 130.677826
$ws.Range("O2").Value = 0.9894183625413969
$ws.Range("P2").Value = 0.9894183625413967
$ws.Range("Q2").Value = 7.169174291219779
$ws.Range("R2").Value = 64.52256862097801
$ws.Range("S2").Value = 0.03448961222152856
$ws.Range("T2").Value = 0.03448961222152856

# Row 3: ECs -> FAPs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Slit2"
$ws.Range("C3").Value = "Robo4"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.1645843333333333
$ws.Range("H3").Value = 0.493753
$ws.Range("I3").Value = 0.03485847193389392
$ws.Range("J3").Value = 0.03485847193389392
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.205596
$ws.Range("N3").Value = 0.616788
$ws.Range("O3").Value = 0.004669968820840217
$ws.Range("P3").Value = 0.004669968820840216
$ws.Range("Q3").Value = 0.033837880596
$ws.Range("R3").Value = 0.304540925364
$ws.Range("S3").Value = 0.0001627879770734184
$ws.Range("T3").Value = 0.0001627879770734183

# Row 4: ECs -> sCs
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Slit2"
$ws.Range("C4").Value = "Robo4"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.1645843333333333
$ws.Range("H4").Value = 0.493753
$ws.Range("I4").Value = 0.03485847193389392
$ws.Range("J4").Value = 0.03485847193389392
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.260262
$ws.Range("N4").Value = 0.780786
$ws.Range("O4").Value = 0.005911668637762975
$ws.Range("P4").Value = 0.005911668637762974
$ws.Range("Q4").Value = 0.042835047762
$ws.Range("R4").Value = 0.385515429858
$ws.Range("S4").Value = 0.0002060717352919415
$ws.Range("T4").Value = 0.0002060717352919415

# Row 5: FAPs -> ECs
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Slit2"
$ws.Range("C5").Value = "Robo4"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 3.368329
$ws.Range("H5").Value = 10.104987
$ws.Range("I5").Value = 0.7134020567608963
$ws.Range("J5").Value = 0.7134020567608964
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 43.55927533333334
$ws.Range("N5").Value = 130.677826
$ws.Range("O5").Value = 0.9894183625413969
$ws.Range("P5").Value = 0.9894183625413967
$ws.Range("Q5").Value = 146.7219703242513
$ws.Range("R5").Value = 1320.497732918262
$ws.Range("S5").Value = 0.7058530948340307
$ws.Range("T5").Value = 0.7058530948340307

# Row 6: FAPs -> FAPs
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Slit2"
$ws.Range("C6").Value = "Robo4"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 3.368329
$ws.Range("H6").Value = 10.104987
$ws.Range("I6").Value = 0.7134020567608963
$ws.Range("J6").Value = 0.7134020567608964
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.205596
$ws.Range("N6").Value = 0.616788
$ws.Range("O6").Value = 0.004669968820840217
$ws.Range("P6").Value = 0.004669968820840216
$ws.Range("Q6").Value = 0.6925149690839999
$ws.Range("R6").Value = 6.232634721756
$ws.Range("S6").Value = 0.003331565361796668
$ws.Range("T6").Value = 0.003331565361796668

# Row 7: FAPs -> sCs
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Slit2"
$ws.Range("C7").Value = "Robo4"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 3.368329
$ws.Range("H7").Value = 10.104987
$ws.Range("I7").Value = 0.7134020567608963
$ws.Range("J7").Value = 0.7134020567608964
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.260262
$ws.Range("N7").Value = 0.780786
$ws.Range("O7").Value = 0.005911668637762975
$ws.Range("P7").Value = 0.005911668637762974
$ws.Range("Q7").Value = 0.8766480421979999
$ws.Range("R7").Value = 7.889832379782
$ws.Range("S7").Value = 0.004217396565068992
$ws.Range("T7").Value = 0.004217396565068992

# Row 8: sCs -> ECs
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Slit2"
$ws.Range("C8").Value = "Robo4"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.188588333333333
$ws.Range("H8").Value = 3.565765
$ws.Range("I8").Value = 0.2517394713052097
$ws.Range("J8").Value = 0.2517394713052098
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 43.55927533333334
$ws.Range("N8").Value = 130.677826
$ws.Range("O8").Value = 0.9894183625413969
$ws.Range("P8").Value = 0.9894183625413967
$ws.Range("Q8").Value = 51.77404646965445
$ws.Range("R8").Value = 465.96641822689
$ws.Range("S8").Value = 0.2490756554858376
$ws.Range("T8").Value = 0.2490756554858376

# Row 9: sCs -> FAPs
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Slit2"
$ws.Range("C9").Value = "Robo4"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.188588333333333
$ws.Range("H9").Value = 3.565765
$ws.Range("I9").Value = 0.2517394713052097
$ws.Range("J9").Value = 0.2517394713052098
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.205596
$ws.Range("N9").Value = 0.616788
$ws.Range("O9").Value = 0.004669968820840217
$ws.Range("P9").Value = 0.004669968820840216
$ws.Range("Q9").Value = 0.24436900698
$ws.Range("R9").Value = 2.19932106282
$ws.Range("S9").Value = 0.00117561548197013
$ws.Range("T9").Value = 0.00117561548197013

# Row 10: sCs -> sCs
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Slit2"
$ws.Range("C10").Value = "Robo4"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.188588333333333
$ws.Range("H10").Value = 3.565765
$ws.Range("I10").Value = 0.2517394713052097
$ws.Range("J10").Value = 0.2517394713052098
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.260262
$ws.Range("N10").Value = 0.780786
$ws.Range("O10").Value = 0.005911668637762975
$ws.Range("P10").Value = 0.005911668637762974
$ws.Range("Q10").Value = 0.30934437681
$ws.Range("R10").Value = 2.78409939129
$ws.Range("S10").Value = 0.001488200337402041
$ws.Range("T10").Value = 0.001488200337402041
